$wb = $excel.ActiveWorkbook

# The sheet holding the "Sygehus PAS-systemer" table is named after the date
# the data was last reviewed/updated. This edit bumps that date forward
# (02-12-2025 -> 05-12-2025); the worksheet rename automatically cascades to
# the workbook's defined name ("Sygehus_PAS_systemer") that references the
# sheet by name.
$ws = $wb.ActiveSheet
$ws.Name = "Opdateret d. 05-12-2025"
